# "updated figure and data"
#
# The old "EU Production" sheet (4 rows of ad-hoc biofuel figures) is
# removed entirely. The sheet formerly named "Sheet1" (13 rows of
# EU domestic-aviation-consumption / biofuel-production data, already
# tab-selected) is renamed to "EU Production" and takes its place right
# after "ReFuelEU". The active selection on that sheet moves from G17 to
# E42.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the old "EU Production" sheet (its data is superseded).
[void]$wb.Worksheets("EU Production").Delete()

# "Sheet1" becomes the new "EU Production" sheet.
$ws = $wb.Worksheets("Sheet1")
$ws.Name = "EU Production"

# Make it the active/selected sheet and update the selected cell.
[void]$ws.Activate()
[void]$ws.Range("E42").Select()
